$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "From excel - " prefix from the Subject value in B3
$ws.Range("B3").Value = "New Employe Boarding Annoucement September 2023"

# Move the active selection to B10 (as seen in the saved file's sheetView)
$ws.Range("B10").Select()
